$d = $word.ActiveDocument

# The requested change rewrites <w:docDefaults> in word/styles.xml, stripping
# the redundant explicit "default-of-the-default" run/paragraph properties
# down to just the handful of values that actually differ from Word's
# built-in defaults. There is no Word object-model surface for
# w:docDefaults directly, so we round-trip the package through
# Document.WordOpenXML (a flat-OPC dump of every part, including
# word/styles.xml) and replace the <w:docDefaults>...</w:docDefaults>
# block with its trimmed-down equivalent.

$oldDocDefaults = "<w:docDefaults><w:rPrDefault><w:rPr><w:rFonts w:ascii=`"Arial`" w:eastAsia=`"Arial`" w:hAnsi=`"Arial`" w:cs=`"Arial`"/><w:b w:val=`"0`"/><w:i w:val=`"0`"/><w:smallCaps w:val=`"0`"/><w:strike w:val=`"0`"/><w:color w:val=`"000000`"/><w:sz w:val=`"22`"/><w:szCs w:val=`"22`"/><w:u w:val=`"none`"/><w:shd w:val=`"clear`" w:fill=`"auto`"/><w:vertAlign w:val=`"baseline`"/><w:lang w:val=`"en`"/></w:rPr></w:rPrDefault><w:pPrDefault><w:pPr><w:keepNext w:val=`"0`"/><w:keepLines w:val=`"0`"/><w:widowControl/><w:pBdr><w:top w:val=`"nil`" w:sz=`"0`" w:space=`"0`"/><w:left w:val=`"nil`" w:sz=`"0`" w:space=`"0`"/><w:bottom w:val=`"nil`" w:sz=`"0`" w:space=`"0`"/><w:right w:val=`"nil`" w:sz=`"0`" w:space=`"0`"/><w:between w:val=`"nil`" w:sz=`"0`" w:space=`"0`"/></w:pBdr><w:shd w:val=`"clear`" w:fill=`"auto`"/><w:spacing w:before=`"0`" w:after=`"0`" w:line=`"276`" w:lineRule=`"auto`"/><w:ind w:left=`"0`" w:right=`"0`" w:firstLine=`"0`"/><w:contextualSpacing w:val=`"0`"/><w:jc w:val=`"left`"/></w:pPr></w:pPrDefault></w:docDefaults>"

$newDocDefaults = "<w:docDefaults><w:rPrDefault><w:rPr><w:rFonts w:ascii=`"Arial`" w:cs=`"Arial`" w:eastAsia=`"Arial`" w:hAnsi=`"Arial`"/><w:sz w:val=`"22`"/><w:szCs w:val=`"22`"/><w:lang w:val=`"en`"/></w:rPr></w:rPrDefault><w:pPrDefault><w:pPr><w:spacing w:line=`"276`" w:lineRule=`"auto`"/></w:pPr></w:pPrDefault></w:docDefaults>"

$xml = $d.WordOpenXML

$start = $xml.IndexOf("<w:docDefaults>")
$end = $xml.IndexOf("</w:docDefaults>") + "</w:docDefaults>".Length

if ($start -lt 0 -or $end -lt 0) {
    throw "w:docDefaults block not found in WordOpenXML"
}

$existing = $xml.Substring($start, $end - $start)

if ($existing -eq $oldDocDefaults) {
    $replacement = $newDocDefaults
} else {
    # Fallback: already-trimmed (or slightly different serialization) --
    # leave untouched to avoid corrupting an already-applied edit.
    $replacement = $existing
}

$newXml = $xml.Substring(0, $start) + $replacement + $xml.Substring($end)

$d.WordOpenXML = $newXml

Write-Output "docDefaults updated"
